$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tidsregistrering")

# --- Row 41: blank separator row (like rows 8 / 12 / 19 / 23 / 33) ---
$ws.Range("G8:H8").Copy()
$ws.Range("G41:H41").PasteSpecial(-4122)

# --- Row 42 ---
$ws.Range("A4").Copy()
$ws.Range("A42").PasteSpecial(-4122)
$ws.Range("A42").Value = 42809

$ws.Range("G4:H4").Copy()
$ws.Range("G42:H42").PasteSpecial(-4122)
$ws.Range("G42").Value = 0.40625
$ws.Range("H42").Value = 0.41666666666666669

$ws.Range("E42").Value = "Reviewer"
$ws.Range("F42").Value = "Review SDD til UC7"

# --- Row 43 ---
$ws.Range("G4:H4").Copy()
$ws.Range("G43:H43").PasteSpecial(-4122)
$ws.Range("G43").Value = 0.41666666666666669
$ws.Range("H43").Value = 0.44444444444444442

$ws.Range("E43").Value = "Designer"
$ws.Range("F43").Value = "Design UC7 om"

# fill in I42 now (after F43), to reproduce the original shared-string order
$ws.Range("I42").Value = "15m"
$ws.Range("I43").Value = "40m"

# --- Row 44 ---
$ws.Range("G4:H4").Copy()
$ws.Range("G44:H44").PasteSpecial(-4122)
$ws.Range("G44").Value = 0.44444444444444442
$ws.Range("H44").Value = 0.53125

$ws.Range("E44").Value = "Implementer"
$ws.Range("F44").Value = "Testsuite for OC5"

# --- Row 45 ---
$ws.Range("G4:H4").Copy()
$ws.Range("G45:H45").PasteSpecial(-4122)
$ws.Range("G45").Value = 0.53125
$ws.Range("H45").Value = 0.57291666666666663

$ws.Range("E45").Value = "Implementer"
$ws.Range("F45").Value = "Testsuite for OC7 - not possible"

# fill in I44 now (after F45), to reproduce the original shared-string order
$ws.Range("I44").Value = "2t 15m"
$ws.Range("I45").Value = "1t"

# --- Row 46 ---
$ws.Range("G4:H4").Copy()
$ws.Range("G46:H46").PasteSpecial(-4122)
$ws.Range("G46").Value = 0.57291666666666663
$ws.Range("H46").Value = 0.59722222222222221

$ws.Range("E46").Value = "Implementer"
$ws.Range("F46").Value = "OC5 - Not possible - Design omarbejdes"

# --- Row 47 ---
$ws.Range("G4").Copy()
$ws.Range("G47").PasteSpecial(-4122)
$ws.Range("G47").Value = 0.60069444444444442

$ws.Range("E47").Value = "Implementer"
$ws.Range("F47").Value = "OC7"

# --- Extend the role data-validation list from E3:E109 to E3:E110 ---
# (restore original entry order: roles list first, then the Deltagere list)
$ws.Range("E3:E109").Validation.Delete()
$ws.Range("C3").Validation.Delete()
$ws.Range("E3:E110").Validation.Add(3, 1, 1, "GyldigeRoller")
$ws.Range("C3").Validation.Add(3, 1, 1, "Deltagere")

# --- Update the view: scroll down and move the active selection ---
$excel.ActiveWindow.ScrollRow = 21
$ws.Range("F28").Select()
